$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# sim_N: simulate more households
$ws.Range("D7").Value = 100

# Equipment presence rates become probabilities instead of booleans
$ws.Range("D20").Value = 0.65   # equipment_WashingMachine
$ws.Range("D21").Value = 0.6    # equipment_TumbleDryer
$ws.Range("D22").Value = 0.6    # equipment_DishWasher
$ws.Range("D23").Value = 0.25   # equipment_WasherDryer

# Slightly taller rows for the equipment block
$ws.Rows.Item(21).RowHeight = 19.5
$ws.Rows.Item(22).RowHeight = 19.5
$ws.Rows.Item(23).RowHeight = 19.5
$ws.Rows.Item(24).RowHeight = 19.5
$ws.Rows.Item(25).RowHeight = 19.5

# EV_present: show one decimal place in its number format
$ws.Range("D30").NumberFormat = "#,##0.0"
